# Applies the "edits some pages and add some wishlist tests" commit:
#  - Duplicates the "productsNames" sheet into a new "productsNames (2)" sheet
#    at the end of the workbook, with new row-2 content, and makes it the
#    active/selected sheet.
#  - The previously active sheet ("contact") is no longer the selected tab.
#  - Updates the workbook view's first visible tab / active tab indices.

$wb = $excel.ActiveWorkbook

# "productsNames" is the template sheet to copy from.
$template = $wb.Worksheets.Item("productsNames")
$template.Copy($wb.Worksheets.Item($wb.Worksheets.Count))

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "productsNames (2)"

# Replace the second row of data with the new product names.
$newSheet.Range("A2").Value = "Small Rainbow Jumper"
$newSheet.Range("B2").Value = "CFrangipani Embroidered Dress"

$newSheet.Range("B2").Select()

# Make the new sheet the active / selected tab, and scroll the tab strip
# so "categoryNames" (index 4, 1-based) is the first visible tab.
$newSheet.Activate()
$wb.Windows.Item(1).ScrollWorkbookTabs(4)
